# Reading Bid opening time in a separate column
#
# Splits the existing "data de abertura" (opening date) column into two
# columns: the date (kept in column E) and a new "Hora de Abertura"
# (opening time) column inserted right after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F: everything that was F:K shifts to G:L.
$ws.Columns.Item(6).Insert()

# New header cell for the inserted column.
$ws.Range("F1").Value = "Hora de Abertura"

# Update the date-only values in column E (the time-of-day portion is
# removed from here and now lives in column F).
$ws.Range("E2").Value = 42479.041666666664
$ws.Range("E3").Value = 42468.583333333336
$ws.Range("E4").Value = 42480.041666666664
$ws.Range("E5").Value = 29609.523611111112

# Populate the new "Hora de Abertura" column with the opening time
# (fraction-of-day) values and a time number format.
$ws.Range("F2").Value = 0.58333333333333337
$ws.Range("F3").Value = 0.5
$ws.Range("F4").Value = 0.41666666666666669
$ws.Range("F5").Value = 0.84930555555555554
$ws.Range("F2:F5").NumberFormat = "h:mm"

# Column widths: E and F both become narrower "date"/"time" columns, and
# the shifted-right J/K/L columns pick up their own explicit widths.
$ws.Columns.Item(5).ColumnWidth = 14.66666666666667
$ws.Columns.Item(6).ColumnWidth = 14.66666666666667
$ws.Columns.Item(10).ColumnWidth = 10.83333333333333
$ws.Columns.Item(11).ColumnWidth = 12.83333333333333
$ws.Columns.Item(12).ColumnWidth = 13.5

# Scroll the view over a bit and select the new time cell for row 3, as
# in the saved workbook.
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollColumn = 5
} catch {
}
$ws.Range("F3").Select()
